# Adds the "Check Header Social Media links are working" test script
# (row 17) test-steps / expected-result detail, mirroring the same detail
# already present for the Footer Social Media test (row 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared-string text for the header social-media test case ---
$testSteps = "1. Open https://abantecart.codifyme.co.nz`n2. I inspect Facebook link`n3. I inspect Twitter link"
$expected  = "2. Link is correct and working (https://www.facebook.com/AbanteCart)`n3. Link is correct and working (https://twitter.com/abantecart)"

# --- Row 17 ("6A" / Check Header Social Media links are working) ---
# Match the formatting already used by rows 2-16 (style index 3): copy it
# from row 16 onto row 17's A:D cells.
$ws.Range("A16:D16").Copy()
$ws.Range("A17:D17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E17").Value = $testSteps
$ws.Range("F17").Value = $expected
$ws.Range("D17").Copy()
$ws.Range("E17").PasteSpecial(-4122)  # xlPasteFormats (keep same style as rest of row)
$ws.Range("D17").Copy()
$ws.Range("F17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(17).RowHeight = 60

# --- Row 18 ("6B" / Check Footer Social Media links are working) ---
# Same two columns of detail, same text, keeping this row's existing style.
$ws.Range("E18").Value = $testSteps
$ws.Range("F18").Value = $expected
$ws.Range("A18").Copy()
$ws.Range("E18:F18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(18).RowHeight = 60

# --- Rows 19, 22 & 25 revert to the sheet's default (auto) row height ---
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(25).AutoFit()

# --- Update the frozen-pane view / current selection ---
$ws.Range("C18").Select() | Out-Null
